$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Covid-19 podatki")

# Expand the "Tabela1" table to include the new row (A1:J80 -> A1:J81)
$lo = $ws.ListObjects.Item("Tabela1")
$lo.Resize($ws.Range("A1:J81"))

# Copy the formatting of the previous row so the new row matches visually
# (date-style format for column A, plain right aligned number format for B..J)
$ws.Range("A80").Copy()
$ws.Range("A81").PasteSpecial(-4122)

$ws.Range("C80:J80").Copy()
$ws.Range("B81:I81").PasteSpecial(-4122)

$ws.Range("J80").Copy()
$ws.Range("J81").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# New day of data (2020-05-31 snapshot referenced in the commit message)
$ws.Range("A81").Value = 43981
$ws.Range("B81").Value = 78793
$ws.Range("C81").Value = 264
$ws.Range("D81").Value = 1473
$ws.Range("E81").Value = 0
$ws.Range("F81").Value = 6
$ws.Range("G81").Value = 2
$ws.Range("H81").Value = 1
$ws.Range("I81").Value = 108
$ws.Range("J81").Value = 0

# Update the view so the new row is visible/selected, mirroring what Excel
# does automatically when a user types data into the row right below a table
$win = $excel.ActiveWindow
$win.ScrollRow = 50
$ws.Range("A81:J81").Select() | Out-Null
